# Add a 4th worksheet "AutoFiltered" that duplicates Sheet1's data, turn on
# an AutoFilter over its data range, and restrict its print area -- while
# keeping the _FilterDatabase / Print_Area defined names scoped locally to
# that sheet only (so GetWorksheetNames() callers don't see Excel's builtin
# table names leak in as if they were worksheet names).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Sheet1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy Sheet1 (identical A1:D8 company data) to the end of the workbook.
$source.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "AutoFiltered"

# Make it the active/selected sheet with the selection Excel leaves behind
# after turning on AutoFilter + clicking back on the sheet.
$newSheet.Activate()
$newSheet.Range("E9").Select()

# Turn on the AutoFilter for the whole data range.
$newSheet.Range("A1:D8").AutoFilter()

# Register the (hidden) _FilterDatabase name AutoFilter implies, scoped to
# this sheet only.
$filterName = $newSheet.Names.Add("_xlnm._FilterDatabase", "=AutoFiltered!`$A`$1:`$D`$8")
$filterName.Visible = $false

# Restrict the print area to the first two columns / six rows, also scoped
# to this sheet only.
$newSheet.PageSetup.PrintArea = "`$A`$1:`$B`$6"
